# Generate Report for Handoff
#
# The localization-status workbook tracks, per target language sheet, the
# "Latest Handoff Datetime" (column D) for every source file that has been
# handed off for translation. A new handoff was generated for the
# "91b88c10-..." file (row 5 of each language sheet), so its previously
# stale/duplicated handoff timestamp needs to be refreshed to the new
# handoff time for both the zh-cn and de-de target languages.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 5 corresponds to the 91b88c10-... source file.
# Record the new "Latest Handoff Datetime" for this handoff.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-04 15:27:18"

# de-de sheet: same source file / row, new handoff was also generated
# for the German target.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-04 15:27:30"
